$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (Buy num) to fit the new longer text value
$ws.Columns.Item(2).ColumnWidth = 19.33

# New note row 23
$ws.Range("A23").Value = "Check R1's value to see if 1k is appropriate"

# Row 12's "Buy num" cell becomes a text note instead of a plain number
$ws.Range("B12").Value = "1 (originally 2)"

# New component row 22
$ws.Range("A22").Value = "300 Ohm capacitor (3.6V)"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = "594-MCT06030C3000FP5"
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = 0.134
$ws.Range("F22").Value = "Thin Film Resistors - SMD .1W 300ohms 1% 0603 50ppm Auto"

# Restore the cursor/selection to where the user last clicked
[void]$ws.Range("F8").Select()
